# Regenerate save_data column G ("K") values, replacing old Strike# derived
# numbers with the newly computed K values (calc and write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 9
    4  = 4
    5  = 5
    6  = 3
    7  = 3
    8  = 6
    9  = 5
    10 = 5
    11 = 6
    12 = 4
    13 = 4
    14 = 4
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
